$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.459.16"
$ws.Range("E2").Value = "  +4.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.838.96"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.029"
$ws.Range("E4").Value = "  +2.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.33"
$ws.Range("E5").Value = "  +4.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.025"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4365"
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3727"
$ws.Range("E8").Value = "  +3.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07372"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8740"
$ws.Range("E10").Value = "  +4.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.40"
$ws.Range("E11").Value = "  +5.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.883.82"
$ws.Range("E12").Value = "  +5.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.471"
$ws.Range("E13").Value = "  +4.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.690"
$ws.Range("E14").Value = "  +3.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07151"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.61"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008992"
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.40"
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.467.46"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.235"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.27"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.085.14"
$ws.Range("E24").Value = "  +3.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.68"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.907"
$ws.Range("E26").Value = "  +5.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.60"
$ws.Range("E27").Value = "  +3.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.244"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.925"
$ws.Range("E29").Value = "  +5.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.06"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09057"
$ws.Range("E31").Value = "  +2.53%  "
$ws.Range("E32").Value = "  +7.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7605"
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.485"
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.866"
$ws.Range("E35").Value = "  +4.73%  "
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.148"
$ws.Range("E37").Value = "  +5.77%  "
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05253"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5171"
$ws.Range("E40").Value = "  +5.31%  "
$ws.Range("E41").Value = "  +6.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1664"
$ws.Range("E42").Value = "  +3.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.550"
$ws.Range("E43").Value = "  +2.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.497"
$ws.Range("E44").Value = "  +6.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.90"
$ws.Range("E45").Value = "  +4.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.56"
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.027"
$ws.Range("E47").Value = "  +2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.685"
$ws.Range("E48").Value = "  +3.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4630"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.898"
$ws.Range("E50").Value = "  +10.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06316"
$ws.Range("E51").Value = "  +2.30%  "
